$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 2440.7334
$ws.Cells.Item(6, 9).Value = 3311.182
$ws.Cells.Item(6, 10).Value = 47
$ws.Cells.Item(6, 11).Value = 9933.545999999998
$ws.Cells.Item(6, 12).Value = 141
$ws.Cells.Item(6, 13).Value = -9821.545999999998
$ws.Cells.Item(6, 14).Value = -365
$ws.Cells.Item(13, 8).Value = 610.25
$ws.Cells.Item(13, 10).Value = 383
$ws.Cells.Item(13, 12).Value = 383
$ws.Cells.Item(13, 14).Value = -721
$ws.Cells.Item(29, 8).Value = 277.25
$ws.Cells.Item(29, 10).Value = 366.33334
$ws.Cells.Item(29, 12).Value = 1099.00002
$ws.Cells.Item(29, 14).Value = -1661.00002
$ws.Cells.Item(32, 8).Value = 4925.625
$ws.Cells.Item(32, 9).Value = 3333
$ws.Cells.Item(32, 10).Value = 5153.143
$ws.Cells.Item(32, 11).Value = 3333
$ws.Cells.Item(32, 12).Value = 5153.143
$ws.Cells.Item(32, 13).Value = -3007
$ws.Cells.Item(32, 14).Value = -5805.143
$ws.Cells.Item(38, 8).Value = 3389.625
$ws.Cells.Item(38, 9).Value = 3588.2144
$ws.Cells.Item(38, 10).Value = 1999.5
$ws.Cells.Item(38, 11).Value = 10764.6432
$ws.Cells.Item(38, 12).Value = 5998.5
$ws.Cells.Item(38, 13).Value = -10392.6432
$ws.Cells.Item(38, 14).Value = -6742.5
$ws.Cells.Item(43, 8).Value = 4999
$ws.Cells.Item(43, 10).Value = 4999
$ws.Cells.Item(43, 12).Value = 4999
$ws.Cells.Item(43, 14).Value = -5137
$ws.Cells.Item(53, 8).Value = 439.1875
$ws.Cells.Item(53, 9).Value = 368.44446
$ws.Cells.Item(53, 10).Value = 530.1429000000001
$ws.Cells.Item(53, 11).Value = 368.44446
$ws.Cells.Item(53, 12).Value = 530.1429000000001
$ws.Cells.Item(53, 13).Value = 268.55554
$ws.Cells.Item(53, 14).Value = -1804.1429
$ws.Cells.Item(62, 8).Value = 6354.2856
$ws.Cells.Item(62, 9).Value = 5907.273
$ws.Cells.Item(62, 10).Value = 7993.3335
$ws.Cells.Item(62, 11).Value = 5907.273
$ws.Cells.Item(62, 12).Value = 7993.3335
$ws.Cells.Item(62, 13).Value = -5283.273
$ws.Cells.Item(62, 14).Value = -9241.333500000001
$ws.Cells.Item(65, 8).Value = 6354.2856
$ws.Cells.Item(65, 9).Value = 5907.273
$ws.Cells.Item(65, 10).Value = 7993.3335
$ws.Cells.Item(65, 11).Value = 29536.365
$ws.Cells.Item(65, 12).Value = 39966.6675
$ws.Cells.Item(65, 13).Value = -26416.365
$ws.Cells.Item(65, 14).Value = -46206.6675
$ws.Cells.Item(86, 8).Value = 4641
$ws.Cells.Item(86, 10).Value = 7249.25
$ws.Cells.Item(86, 12).Value = 7249.25
$ws.Cells.Item(86, 14).Value = -9495.25
$ws.Cells.Item(89, 8).Value = 4641
$ws.Cells.Item(89, 10).Value = 7249.25
$ws.Cells.Item(89, 12).Value = 36246.25
$ws.Cells.Item(89, 14).Value = -47478.25
$ws.Cells.Item(96, 8).Value = 16318.23
$ws.Cells.Item(96, 9).Value = 22716.445
$ws.Cells.Item(96, 11).Value = 68149.33499999999
$ws.Cells.Item(96, 13).Value = -66776.33499999999
$ws.Cells.Item(107, 8).Value = 409.5263
$ws.Cells.Item(107, 9).Value = 400.13333
$ws.Cells.Item(107, 10).Value = 444.75
$ws.Cells.Item(107, 11).Value = 400.13333
$ws.Cells.Item(107, 12).Value = 444.75
$ws.Cells.Item(107, 13).Value = 1519.86667
$ws.Cells.Item(107, 14).Value = -4284.75
$ws.Cells.Item(125, 8).Value = 187502130
$ws.Cells.Item(125, 9).Value = 200002000
$ws.Cells.Item(125, 11).Value = 1800018000
$ws.Cells.Item(125, 13).Value = -1800015540
$ws.Cells.Item(132, 8).Value = 1465.9333
$ws.Cells.Item(132, 9).Value = 1323.3846
$ws.Cells.Item(132, 11).Value = 3970.1538
$ws.Cells.Item(132, 13).Value = -1440.1538
$ws.Cells.Item(135, 8).Value = 1780.2941
$ws.Cells.Item(135, 10).Value = 4997.5
$ws.Cells.Item(135, 12).Value = 44977.5
$ws.Cells.Item(135, 14).Value = -50047.5
$ws.Cells.Item(137, 8).Value = 2127.5454
$ws.Cells.Item(137, 9).Value = 1940.1
$ws.Cells.Item(137, 10).Value = 4002
$ws.Cells.Item(137, 11).Value = 5820.299999999999
$ws.Cells.Item(137, 12).Value = 12006
$ws.Cells.Item(137, 13).Value = -3270.299999999999
$ws.Cells.Item(137, 14).Value = -17106
$ws.Cells.Item(138, 8).Value = 4599.13
$ws.Cells.Item(138, 9).Value = 3521.2856
$ws.Cells.Item(138, 10).Value = 5003.3213
$ws.Cells.Item(138, 11).Value = 10563.8568
$ws.Cells.Item(138, 12).Value = 15009.9639
$ws.Cells.Item(138, 13).Value = -5423.856800000001
$ws.Cells.Item(138, 14).Value = -25289.9639

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 9400.275
$ws.Cells.Item(32, 9).Value = 6908.385
$ws.Cells.Item(32, 10).Value = 30996.666
$ws.Cells.Item(32, 11).Value = 6908.385
$ws.Cells.Item(32, 12).Value = 30996.666
$ws.Cells.Item(32, 13).Value = -6621.385
$ws.Cells.Item(32, 14).Value = -31570.666
$ws.Cells.Item(74, 8).Value = 2522.075
$ws.Cells.Item(74, 9).Value = 2390.7715
$ws.Cells.Item(74, 11).Value = 2390.7715
$ws.Cells.Item(74, 13).Value = -1516.7715
$ws.Cells.Item(77, 8).Value = 2522.075
$ws.Cells.Item(77, 9).Value = 2390.7715
$ws.Cells.Item(77, 11).Value = 11953.8575
$ws.Cells.Item(77, 13).Value = -7585.857499999998
$ws.Cells.Item(102, 8).Value = 1781.5385
$ws.Cells.Item(102, 9).Value = 1781.5385
$ws.Cells.Item(102, 10).Value = 0
$ws.Cells.Item(102, 11).Value = 1781.5385
$ws.Cells.Item(102, 12).Value = 0
$ws.Cells.Item(102, 13).ClearContents()
$ws.Cells.Item(102, 14).Value = -159.5385000000001
$ws.Cells.Item(110, 8).Value = 630.35297
$ws.Cells.Item(110, 9).Value = 630.35297
$ws.Cells.Item(110, 11).Value = 630.35297
$ws.Cells.Item(110, 13).Value = 1414.64703
$ws.Cells.Item(122, 8).Value = 2664.125
$ws.Cells.Item(122, 9).Value = 2285.5
$ws.Cells.Item(122, 11).Value = 6856.5
$ws.Cells.Item(122, 13).Value = -4406.5
$ws.Cells.Item(127, 8).Value = 40000
$ws.Cells.Item(127, 9).Value = 40000
$ws.Cells.Item(127, 11).Value = 40000
$ws.Cells.Item(127, 13).Value = -35040
$ws.Cells.Item(132, 8).Value = 2443.8
$ws.Cells.Item(132, 9).Value = 2388.2104
$ws.Cells.Item(132, 10).Value = 3500
$ws.Cells.Item(132, 11).Value = 7164.6312
$ws.Cells.Item(132, 12).Value = 10500
$ws.Cells.Item(132, 13).Value = -4634.6312
$ws.Cells.Item(132, 14).Value = -15560
$ws.Cells.Item(141, 8).Value = 99999
$ws.Cells.Item(141, 10).Value = 99999
$ws.Cells.Item(141, 12).Value = 99999
$ws.Cells.Item(141, 14).Value = -110359

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 2395.9375
$ws.Cells.Item(99, 9).Value = 1702.4166
$ws.Cells.Item(99, 11).Value = 1702.4166
$ws.Cells.Item(99, 13).Value = -204.4166
$ws.Cells.Item(103, 8).Value = 36051.145
$ws.Cells.Item(103, 10).Value = 36051.145
$ws.Cells.Item(103, 12).Value = 36051.145
$ws.Cells.Item(103, 14).Value = -38395.145
$ws.Cells.Item(105, 8).Value = 2245.6667
$ws.Cells.Item(105, 9).Value = 1694.8
$ws.Cells.Item(105, 11).Value = 1694.8
$ws.Cells.Item(105, 13).Value = 52.20000000000005
$ws.Cells.Item(107, 8).Value = 1818.3
$ws.Cells.Item(107, 10).Value = 3353
$ws.Cells.Item(107, 12).Value = 3353
$ws.Cells.Item(107, 14).Value = -7193
$ws.Cells.Item(134, 8).Value = 3824.5
$ws.Cells.Item(134, 9).Value = 3377.4285
$ws.Cells.Item(134, 11).Value = 10132.2855
$ws.Cells.Item(134, 13).Value = -7597.2855
$ws.Cells.Item(135, 8).Value = 56598.6
$ws.Cells.Item(135, 10).Value = 74331.664
$ws.Cells.Item(135, 12).Value = 74331.664
$ws.Cells.Item(135, 14).Value = -84471.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 46.285713
$ws.Cells.Item(7, 9).Value = 48.117645
$ws.Cells.Item(7, 10).Value = 38.5
$ws.Cells.Item(7, 11).Value = 48.117645
$ws.Cells.Item(7, 12).Value = 38.5
$ws.Cells.Item(7, 13).Value = 64.88235499999999
$ws.Cells.Item(7, 14).Value = -264.5
$ws.Cells.Item(16, 8).Value = 2847.6
$ws.Cells.Item(16, 9).Value = 3071
$ws.Cells.Item(16, 10).Value = 2326.3333
$ws.Cells.Item(16, 11).Value = 3071
$ws.Cells.Item(16, 12).Value = 2326.3333
$ws.Cells.Item(16, 13).Value = -2784
$ws.Cells.Item(16, 14).Value = -2900.3333
$ws.Cells.Item(22, 8).Value = 10466.36
$ws.Cells.Item(22, 9).Value = 334.29413
$ws.Cells.Item(22, 11).Value = 334.29413
$ws.Cells.Item(22, 13).Value = 15.70587
$ws.Cells.Item(31, 8).Value = 17891.076
$ws.Cells.Item(31, 9).Value = 29260.8
$ws.Cells.Item(31, 10).Value = 10785
$ws.Cells.Item(31, 11).Value = 29260.8
$ws.Cells.Item(31, 12).Value = 10785
$ws.Cells.Item(31, 13).Value = -28965.8
$ws.Cells.Item(31, 14).Value = -11375
$ws.Cells.Item(34, 8).Value = 17891.076
$ws.Cells.Item(34, 9).Value = 29260.8
$ws.Cells.Item(34, 10).Value = 10785
$ws.Cells.Item(34, 11).Value = 29260.8
$ws.Cells.Item(34, 12).Value = 10785
$ws.Cells.Item(34, 13).Value = -29058.8
$ws.Cells.Item(34, 14).Value = -11189
$ws.Cells.Item(58, 8).Value = 0
$ws.Cells.Item(58, 9).Value = 0
$ws.Cells.Item(58, 11).Value = 0
$ws.Cells.Item(58, 13).ClearContents()
$ws.Cells.Item(74, 8).Value = 41131.2
$ws.Cells.Item(74, 10).Value = 41131.2
$ws.Cells.Item(74, 12).Value = 41131.2
$ws.Cells.Item(74, 14).Value = -42879.2
$ws.Cells.Item(77, 8).Value = 41131.2
$ws.Cells.Item(77, 10).Value = 41131.2
$ws.Cells.Item(77, 12).Value = 123393.6
$ws.Cells.Item(77, 14).Value = -132129.6
$ws.Cells.Item(105, 8).Value = 3351.6365
$ws.Cells.Item(105, 9).Value = 2588.6
$ws.Cells.Item(105, 10).Value = 3987.5
$ws.Cells.Item(105, 11).Value = 2588.6
$ws.Cells.Item(105, 12).Value = 3987.5
$ws.Cells.Item(105, 13).Value = -841.5999999999999
$ws.Cells.Item(105, 14).Value = -7481.5
$ws.Cells.Item(113, 8).Value = 2847.6
$ws.Cells.Item(113, 9).Value = 3071
$ws.Cells.Item(113, 10).Value = 2326.3333
$ws.Cells.Item(113, 11).Value = 3071
$ws.Cells.Item(113, 12).Value = 2326.3333
$ws.Cells.Item(113, 13).Value = -901
$ws.Cells.Item(113, 14).Value = -6666.3333
$ws.Cells.Item(120, 8).Value = 48999
$ws.Cells.Item(120, 10).Value = 48999
$ws.Cells.Item(120, 12).Value = 48999
$ws.Cells.Item(120, 14).Value = -56257
$ws.Cells.Item(132, 8).Value = 2637
$ws.Cells.Item(132, 9).Value = 2637
$ws.Cells.Item(132, 11).Value = 7911
$ws.Cells.Item(132, 13).Value = -5381
$ws.Cells.Item(134, 8).Value = 2255.1177
$ws.Cells.Item(134, 9).Value = 2158.5625
$ws.Cells.Item(134, 11).Value = 6475.6875
$ws.Cells.Item(134, 13).Value = -3940.6875
$ws.Cells.Item(136, 8).Value = 0
$ws.Cells.Item(136, 9).Value = 0
$ws.Cells.Item(136, 11).Value = 0
$ws.Cells.Item(136, 13).ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 130071.18
$ws.Cells.Item(2, 10).Value = 158698
$ws.Cells.Item(2, 12).Value = 952188
$ws.Cells.Item(2, 14).Value = -952414
$ws.Cells.Item(8, 8).Value = 983.6667
$ws.Cells.Item(8, 9).Value = 983.6667
$ws.Cells.Item(8, 11).Value = 2951.0001
$ws.Cells.Item(8, 13).Value = -2812.0001
$ws.Cells.Item(34, 8).Value = 3998
$ws.Cells.Item(34, 10).Value = 3998
$ws.Cells.Item(34, 12).Value = 11994
$ws.Cells.Item(34, 14).Value = -12162
$ws.Cells.Item(39, 8).Value = 7129.3
$ws.Cells.Item(39, 10).Value = 6943.6665
$ws.Cells.Item(39, 12).Value = 20830.9995
$ws.Cells.Item(39, 14).Value = -21418.9995
$ws.Cells.Item(55, 8).Value = 2640
$ws.Cells.Item(55, 10).Value = 4500
$ws.Cells.Item(55, 12).Value = 13500
$ws.Cells.Item(55, 14).Value = -13854
$ws.Cells.Item(114, 8).Value = 2822.6667
$ws.Cells.Item(114, 9).Value = 2984.4
$ws.Cells.Item(114, 10).Value = 2707.1428
$ws.Cells.Item(114, 11).Value = 8953.200000000001
$ws.Cells.Item(114, 12).Value = 8121.428400000001
$ws.Cells.Item(114, 13).Value = -5699.200000000001
$ws.Cells.Item(114, 14).Value = -14629.4284
$ws.Cells.Item(122, 8).Value = 591.25
$ws.Cells.Item(122, 9).Value = 561.4286
$ws.Cells.Item(122, 10).Value = 800
$ws.Cells.Item(122, 11).Value = 5052.8574
$ws.Cells.Item(122, 12).Value = 7200
$ws.Cells.Item(122, 13).Value = -2602.8574
$ws.Cells.Item(122, 14).Value = -12100
$ws.Cells.Item(129, 8).Value = 2377.375
$ws.Cells.Item(129, 9).Value = 753.1667
$ws.Cells.Item(129, 11).Value = 2259.5001
$ws.Cells.Item(129, 13).Value = 2740.4999
$ws.Cells.Item(132, 8).Value = 5999.7144
$ws.Cells.Item(132, 9).Value = 6799.6
$ws.Cells.Item(132, 10).Value = 4000
$ws.Cells.Item(132, 11).Value = 61196.4
$ws.Cells.Item(132, 12).Value = 36000
$ws.Cells.Item(132, 13).Value = -58666.4
$ws.Cells.Item(132, 14).Value = -41060

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(40, 8).Value = 499
$ws.Cells.Item(40, 10).Value = 499
$ws.Cells.Item(40, 12).Value = 499
$ws.Cells.Item(40, 14).Value = -801
$ws.Cells.Item(48, 8).Value = 34985
$ws.Cells.Item(48, 10).Value = 34985
$ws.Cells.Item(48, 12).Value = 34985
$ws.Cells.Item(48, 14).Value = -35955
$ws.Cells.Item(80, 8).Value = 7729.8096
$ws.Cells.Item(80, 9).Value = 1397
$ws.Cells.Item(80, 11).Value = 1397
$ws.Cells.Item(80, 13).Value = -399
$ws.Cells.Item(83, 8).Value = 7729.8096
$ws.Cells.Item(83, 9).Value = 1397
$ws.Cells.Item(83, 11).Value = 6985
$ws.Cells.Item(83, 13).Value = -1993
$ws.Cells.Item(97, 8).Value = 440.84616
$ws.Cells.Item(97, 9).Value = 403.0909
$ws.Cells.Item(97, 11).Value = 403.0909
$ws.Cells.Item(97, 13).Value = 92.90910000000002
$ws.Cells.Item(102, 8).Value = 2407.25
$ws.Cells.Item(102, 9).Value = 2198.3333
$ws.Cells.Item(102, 11).Value = 2198.3333
$ws.Cells.Item(102, 13).Value = -576.3332999999998
$ws.Cells.Item(107, 8).Value = 481.16666
$ws.Cells.Item(107, 9).Value = 528
$ws.Cells.Item(107, 11).Value = 528
$ws.Cells.Item(107, 13).Value = 1392
$ws.Cells.Item(111, 8).Value = 150000
$ws.Cells.Item(111, 10).Value = 150000
$ws.Cells.Item(111, 12).Value = 150000
$ws.Cells.Item(111, 14).Value = -156134
$ws.Cells.Item(122, 8).Value = 3145.3333
$ws.Cells.Item(122, 9).Value = 3443.5
$ws.Cells.Item(122, 11).Value = 10330.5
$ws.Cells.Item(122, 13).Value = -7880.5
$ws.Cells.Item(126, 8).Value = 2972.111
$ws.Cells.Item(126, 9).Value = 2843.625
$ws.Cells.Item(126, 11).Value = 8530.875
$ws.Cells.Item(126, 13).Value = -6060.875
$ws.Cells.Item(132, 8).Value = 1640.3334
$ws.Cells.Item(132, 9).Value = 1516.7273
$ws.Cells.Item(132, 10).Value = 3000
$ws.Cells.Item(132, 11).Value = 4550.1819
$ws.Cells.Item(132, 12).Value = 9000
$ws.Cells.Item(132, 13).Value = -2020.1819
$ws.Cells.Item(132, 14).Value = -14060

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 3133.125
$ws.Cells.Item(22, 9).Value = 3182.5
$ws.Cells.Item(22, 10).Value = 3083.75
$ws.Cells.Item(22, 11).Value = 3182.5
$ws.Cells.Item(22, 12).Value = 3083.75
$ws.Cells.Item(22, 13).Value = -2887.5
$ws.Cells.Item(22, 14).Value = -3673.75
$ws.Cells.Item(27, 8).Value = 3133.125
$ws.Cells.Item(27, 9).Value = 3182.5
$ws.Cells.Item(27, 10).Value = 3083.75
$ws.Cells.Item(27, 11).Value = 3182.5
$ws.Cells.Item(27, 12).Value = 3083.75
$ws.Cells.Item(27, 13).Value = -3075.5
$ws.Cells.Item(27, 14).Value = -3297.75
$ws.Cells.Item(55, 8).Value = 1622.3334
$ws.Cells.Item(55, 10).Value = 1708.75
$ws.Cells.Item(55, 12).Value = 1708.75
$ws.Cells.Item(55, 14).Value = -2054.75
$ws.Cells.Item(61, 8).Value = 5581.1577
$ws.Cells.Item(61, 9).Value = 5077.625
$ws.Cells.Item(61, 10).Value = 8266.666999999999
$ws.Cells.Item(61, 11).Value = 5077.625
$ws.Cells.Item(61, 12).Value = 8266.666999999999
$ws.Cells.Item(61, 13).Value = -4875.625
$ws.Cells.Item(61, 14).Value = -8670.666999999999
$ws.Cells.Item(82, 8).Value = 1416.5333
$ws.Cells.Item(82, 9).Value = 1025.2
$ws.Cells.Item(82, 10).Value = 2199.2
$ws.Cells.Item(82, 11).Value = 1025.2
$ws.Cells.Item(82, 12).Value = 2199.2
$ws.Cells.Item(82, 13).Value = -664.2
$ws.Cells.Item(82, 14).Value = -2921.2
$ws.Cells.Item(85, 8).Value = 1416.5333
$ws.Cells.Item(85, 9).Value = 1025.2
$ws.Cells.Item(85, 10).Value = 2199.2
$ws.Cells.Item(85, 11).Value = 1025.2
$ws.Cells.Item(85, 12).Value = 2199.2
$ws.Cells.Item(85, 13).Value = 222.8
$ws.Cells.Item(85, 14).Value = -4695.2
$ws.Cells.Item(113, 8).Value = 5581.1577
$ws.Cells.Item(113, 9).Value = 5077.625
$ws.Cells.Item(113, 10).Value = 8266.666999999999
$ws.Cells.Item(113, 11).Value = 5077.625
$ws.Cells.Item(113, 12).Value = 8266.666999999999
$ws.Cells.Item(113, 13).Value = -2907.625
$ws.Cells.Item(113, 14).Value = -12606.667
$ws.Cells.Item(132, 8).Value = 3625.6775
$ws.Cells.Item(132, 9).Value = 3374.9285
$ws.Cells.Item(132, 11).Value = 10124.7855
$ws.Cells.Item(132, 13).Value = -7594.7855
$ws.Cells.Item(133, 8).Value = 96248.664
$ws.Cells.Item(133, 9).Value = 96248.664
$ws.Cells.Item(133, 11).Value = 96248.664
$ws.Cells.Item(133, 13).Value = -93718.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 16667866
$ws.Cells.Item(100, 9).Value = 33333984
$ws.Cells.Item(100, 10).Value = 1749.6666
$ws.Cells.Item(100, 11).Value = 66667968
$ws.Cells.Item(100, 12).Value = 3499.3332
$ws.Cells.Item(100, 13).Value = -66667427
$ws.Cells.Item(100, 14).Value = -4581.3332
$ws.Cells.Item(123, 8).Value = 99995
$ws.Cells.Item(123, 10).Value = 99995
$ws.Cells.Item(123, 12).Value = 99995
$ws.Cells.Item(123, 14).Value = -109795
$ws.Cells.Item(125, 8).Value = 150000
$ws.Cells.Item(125, 10).Value = 150000
$ws.Cells.Item(125, 12).Value = 150000
$ws.Cells.Item(125, 14).Value = -159840
$ws.Cells.Item(132, 8).Value = 8457.5
$ws.Cells.Item(132, 9).Value = 8457.5
$ws.Cells.Item(132, 11).Value = 25372.5
$ws.Cells.Item(132, 13).Value = -22842.5
